$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dates = @("29-10-2021", "30-10-2021", "31-10-2021", "01-11-2021", "02-11-2021")

$startRow = 276
for ($i = 0; $i -lt $dates.Length; $i++) {
    $r = $startRow + $i
    $cellA = $ws.Cells.Item($r, 1)
    # Write the date as a text-returning formula, then paste-special as
    # values. This lands a plain text cell (no number/date auto-detect,
    # no NumberFormat / quote-prefix style side effects) matching how the
    # original date-label column was authored.
    $cellA.Formula = '="' + $dates[$i] + '"'
    $cellA.Copy()
    $cellA.PasteSpecial(-4163)
    $ws.Cells.Item($r, 2).Value = 449
    $ws.Cells.Item($r, 3).Value = 0
}
$excel.CutCopyMode = 0
